# Applies the two kinds of changes described by the commit:
#   1. The "datetimeFigureOut" date placeholder text (Slide Master + all
#      Slide Layouts) is bumped from 2020-10-08 to 2020-11-16.
#   2. Five rectangles on slide 1 get new solid fill colors.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" field text everywhere it is
#    used as a Date Placeholder: the Slide Master and every Custom
#    Layout attached to it.
# ---------------------------------------------------------------------
$oldDate = "2020-10-08"
$newDate = "2020-11-16"

$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Recolor the five stacked rectangles on slide 1 (the pale "taxa"
#    legend blocks). Colors given as BGR-packed ints matching the new
#    srgbClr hex values (D39FE5, 89D1E7, 94D094, F4C37C, EEA0A0).
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

$colorMap = @{
    1 = 15048659   # Rectangle 207: FFA2A2 -> D39FE5
    2 = 15192457   # Rectangle 205: F0B682 -> 89D1E7
    3 = 9752724    # Rectangle 203: C6CE86 -> 94D094
    4 = 8176628    # Rectangle 3  : 90DFB0 -> F4C37C
    5 = 10526958   # Rectangle 1  : 62E8EC -> EEA0A0
}

foreach ($idx in $colorMap.Keys) {
    $shp = $s.Shapes.Item($idx)
    $shp.Fill.ForeColor.RGB = $colorMap[$idx]
}
